$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that gets bumped by
# one day for every row (2 through 486) on each automated run.
$ws.Range("C2:C486").Value = 46061
